# Add a "Units" column to the Concentrations sheet (wc_lang.Concentration
# gained a new `units` attribute - every concentration is now annotated
# with its unit, "M" for molar, in this fixture).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concentrations")

# Insert a new column in front of the existing "Comments" column (C) so the
# layout becomes: Species | Value | Units | Comments | References
$ws.Columns.Item(3).Insert()

$ws.Cells.Item(1, 3).Value = "Units"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = "M"
}

# Leave the editor focused on the sheet/cell that was just edited.
$ws.Select() | Out-Null
$ws.Range("E5").Select() | Out-Null
